$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking values that are stored as TEXT
# (shared strings) in the workbook. Force text formatting before
# assigning so Excel keeps them as text instead of auto-converting
# them to numbers.

$updates = @{
    "C11" = "18.19"
    "B30" = "71.48"
    "C30" = "18.48"
    "D30" = "89.96"
    "B31" = "16.49"
    "C31" = "41.02"
    "D31" = "57.51"
    "B33" = "79.41"
    "C33" = "20.53"
    "D33" = "99.94"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
